$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "BTC", "Bitcoin", 41903, 819830538014, 11986194231, -0.93288)
    ,@(3, "ETH", "Ethereum", 2215.2, 265963134104, 7828168667, -1.48965)
    ,@(4, "USDT", "Tether", 0.999436, 90850654807, 28581861881, 0.0507)
    ,@(5, "BNB", "BNB", 241.3, 37104049349, 540505107, -2.17504)
    ,@(6, "XRP", "XRP", 0.618808, 33380399594, 744026594, -0.72162)
    ,@(7, "SOL", "Solana", 72.72, 31082596893, 1609319506, -3.14503)
    ,@(8, "USDC", "USDC", 1, 24629116934, 4662777447, -0.04887)
    ,@(9, "ADA", "Cardano", 0.598587, 20918001986, 521318006, -2.58032)
    ,@(10, "STETH", "Lido Staked Ether", 2212.43, 20363989148, 33663356, -1.58137)
    ,@(11, "AVAX", "Avalanche", 41.95, 15378544016, 1663507734, -3.43542)
    ,@(12, "DOGE", "Dogecoin", 0.095038, 13491081604, 718179527, 0.87824)
    ,@(13, "DOT", "Polkadot", 6.98, 9113772784, 290157745, -2.58047)
    ,@(14, "TRX", "TRON", 0.103019, 9096921069, 316781809, 0.32609)
    ,@(15, "LINK", "Chainlink", 14.19, 7888449740, 397109911, -2.43447)
    ,@(16, "MATIC", "Polygon", 0.832436, 7715664850, 296001104, -2.53698)
    ,@(17, "WBTC", "Wrapped Bitcoin", 41832, 6469794031, 76901311, -1.04351)
    ,@(18, "SHIB", "Shiba Inu", 0.00001073, 6284307348, 842700935, 7.09315)
    ,@(19, "LTC", "Litecoin", 72.55, 5356673447, 465843609, 0.43497)
    ,@(20, "DAI", "Dai", 0.997801, 5271229056, 199752100, -0.07549)
    ,@(21, "ICP", "Internet Computer", 11.05, 4894968201, 708763029, 20.67143)
    ,@(22, "UNI", "Uniswap", 6.2, 4662259783, 164051700, 0.39648)
    ,@(23, "BCH", "Bitcoin Cash", 228.55, 4473661291, 116735430, -0.80994)
    ,@(24, "LEO", "LEO Token", 3.79, 3513093157, 924406, 0.47316)
    ,@(25, "XLM", "Stellar", 0.12429, 3501286906, 61499931, -0.53519)
    ,@(26, "OKB", "OKB", 57.27, 3437124414, 10572128, -0.88272)
    ,@(27, "ATOM", "Cosmos Hub", 11.62, 3412192174, 285810926, 1.79435)
    ,@(28, "XMR", "Monero", 169.42, 3070847929, 64575945, -1.35297)
    ,@(29, "ETC", "Ethereum Classic", 20.21, 2884407004, 94711177, -0.78112)
    ,@(30, "FIL", "Filecoin", 5.68, 2761809364, 851350434, 5.75791)
    ,@(31, "HBAR", "Hedera", 0.079457, 2660696932, 65957945, -3.12502)
    ,@(32, "CRO", "Cronos", 0.099772, 2636325655, 11076988, -2.29744)
    ,@(33, "IMX", "Immutable", 2.03, 2613492214, 158354376, -6.52462)
    ,@(34, "INJ", "Injective", 30.14, 2530781982, 251068783, -2.50865)
    ,@(35, "APT", "Aptos", 8.17, 2483820305, 139336047, -1.90593)
    ,@(36, "TUSD", "TrueUSD", 0.996213, 2431725952, 127083746, -0.24066)
    ,@(37, "KAS", "Kaspa", 0.109192, 2382753721, 47991154, -10.27094)
    ,@(38, "NEAR", "NEAR Protocol", 2.37, 2382299533, 215167572, 2.44583)
    ,@(39, "VET", "VeChain", 0.02983601, 2166426628, 55024016, -3.7473)
    ,@(40, "TIA", "Celestia", 13.65, 1986988583, 319276355, -0.93387)
    ,@(41, "OP", "Optimism", 2.12, 1928913091, 222874379, -1.929)
    ,@(42, "LDO", "Lido DAO", 2.12, 1879284878, 46221938, -3.26419)
    ,@(43, "TAO", "Bittensor", 316.42, 1856107096, 6125048, -7.10268)
    ,@(44, "MNT", "Mantle", 0.573271, 1793498467, 27516972, -2.00601)
    ,@(45, "RUNE", "THORChain", 5.63, 1690565929, 243385448, -2.94588)
    ,@(46, "EGLD", "MultiversX", 63.25, 1660101565, 49650886, 2.50615)
    ,@(47, "FDUSD", "First Digital USD", 0.998126, 1650714700, 2190988124, -0.11755)
    ,@(48, "QNT", "Quant", 111.9, 1625742004, 18684808, -1.28326)
    ,@(49, "RNDR", "Render", 4.25, 1586314388, 124761442, -4.54793)
    ,@(50, "ALGO", "Algorand", 0.19567, 1564042085, 66958706, -3.32793)
    ,@(51, "AAVE", "Aave", 104.45, 1521640124, 156756034, -3.21341)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
